# This script updates the "想去人数" (want-to-go count) values in column F
# across the four worksheets of the 广州-漫展信息 workbook, reflecting a
# refreshed scrape of the source data (gh-pages output regenerated at 456a3b4).
#
# Sheets:
#   展览     (Exhibitions)
#   演出     (Performances)
#   本地生活 (Local life)
#   全部类型 (All types - aggregate of the above three)

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 26096
$ws.Range("F6").Value = 246
$ws.Range("F7").Value = 584
$ws.Range("F8").Value = 169
$ws.Range("F12").Value = 212
$ws.Range("F15").Value = 285
$ws.Range("F16").Value = 36
$ws.Range("F17").Value = 353
$ws.Range("F19").Value = 1492
$ws.Range("F20").Value = 168
$ws.Range("F21").Value = 16
$ws.Range("F22").Value = 417
$ws.Range("F23").Value = 95

# ---- Sheet: 演出 (Performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 4499
$ws.Range("F3").Value = 222
$ws.Range("F4").Value = 4
$ws.Range("F6").Value = 124
$ws.Range("F15").Value = 45

# ---- Sheet: 本地生活 (Local life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 4940
$ws.Range("F3").Value = 193

# ---- Sheet: 全部类型 (All types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 4940
$ws.Range("F5").Value = 193
$ws.Range("F6").Value = 26096
$ws.Range("F8").Value = 4499
$ws.Range("F9").Value = 246
$ws.Range("F10").Value = 222
$ws.Range("F11").Value = 584
$ws.Range("F12").Value = 4
$ws.Range("F14").Value = 169
$ws.Range("F15").Value = 124
$ws.Range("F16").Value = 124
$ws.Range("F25").Value = 212
$ws.Range("F29").Value = 285
$ws.Range("F30").Value = 36
$ws.Range("F33").Value = 353
$ws.Range("F35").Value = 45
$ws.Range("F36").Value = 1492
$ws.Range("F37").Value = 168
$ws.Range("F39").Value = 16
$ws.Range("F40").Value = 417
$ws.Range("F41").Value = 95
